# Daily attendance processing - 2026-02-18 11:48:34 UTC
# Refresh the session-analysis report: recompute "Recorded By" -> academic
# year coverage per session, refresh class/group statistics, and promote
# two previously-unrecorded sessions (row 21 "C1" / row 49 "C2") to
# Recorded now that attendance records have landed for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Column widths - G narrower now that it holds academic years instead
#    of long recorder-name lists; I narrower to match its shorter values.
# ---------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 42
$ws.Columns.Item(9).ColumnWidth = 10

# ---------------------------------------------------------------------
# 2) Row 21 (Year 3 / C1 / PHARMACOLOGY / session 2) moves from
#    "Not Recorded" (pink highlight) to "Recorded" (normal style),
#    picking up its first attendance record.
# ---------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A21:I21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H21").Value = "1/221"
$ws.Range("I21").Value = "Recorded"

# ---------------------------------------------------------------------
# 3) Row 49 (Year 3 / C2 / PHARMACOLOGY / session 2) gets the same
#    promotion, now with 5 attendance records in.
# ---------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A49:I49").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H49").Value = "5/246"
$ws.Range("I49").Value = "Recorded"

# ---------------------------------------------------------------------
# 4) "Recorded By" column: replace the (now stale) recorder-name lists
#    with the academic years that have contributed attendance data for
#    each session.
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "2025/2026"
$ws.Range("G3").Value = "2022/2023, 2025/2026"
$ws.Range("G4").Value = "2025/2026"
$ws.Range("G5").Value = "2025/2026"
$ws.Range("G6").Value = "2025/2026"
$ws.Range("G7").Value = "2025/2026"
$ws.Range("G8").Value = "2025/2026"
$ws.Range("G9").Value = "2025/2026"
$ws.Range("G10").Value = "2025/2026"
$ws.Range("G11").Value = "2025/2026"
$ws.Range("G12").Value = "2025/2026"
$ws.Range("G13").Value = "2025/2026"
$ws.Range("G14").Value = "2025/2026"
$ws.Range("G15").Value = "2025/2026"
$ws.Range("G16").Value = "2025/2026"
$ws.Range("G17").Value = "2025/2026"
$ws.Range("G18").Value = "2025/2026"
$ws.Range("G19").Value = "2025/2026"
$ws.Range("G20").Value = "2025/2026"
$ws.Range("G21").Value = "2025/2026"
$ws.Range("G22").Value = "2024/2025, 2025/2026"
$ws.Range("G23").Value = "2022/2023, 2023/2024, 2025/2026"
$ws.Range("G24").Value = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
$ws.Range("G25").Value = "2025/2026"
$ws.Range("G26").Value = "2025/2026"
$ws.Range("G27").Value = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
$ws.Range("G28").Value = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
$ws.Range("G29").Value = "2025/2026"
$ws.Range("G30").Value = "2025/2026"
$ws.Range("G31").Value = "2022/2023, 2025/2026"
$ws.Range("G32").Value = "2025/2026"
$ws.Range("G33").Value = "2025/2026"
$ws.Range("G34").Value = "2025/2026"
$ws.Range("G35").Value = "2025/2026"
$ws.Range("G36").Value = "2025/2026"
$ws.Range("G37").Value = "2025/2026"
$ws.Range("G38").Value = "2025/2026"
$ws.Range("G39").Value = "2025/2026"
$ws.Range("G40").Value = "2025/2026"
$ws.Range("G41").Value = "2025/2026"
$ws.Range("G42").Value = "2025/2026"
$ws.Range("G43").Value = "2025/2026"
$ws.Range("G44").Value = "2025/2026"
$ws.Range("G45").Value = "2025/2026"
$ws.Range("G46").Value = "2025/2026"
$ws.Range("G47").Value = "2025/2026"
$ws.Range("G48").Value = "2025/2026"
$ws.Range("G49").Value = "2025/2026"
$ws.Range("G50").Value = "2024/2025, 2025/2026"
$ws.Range("G51").Value = "2022/2023, 2023/2024, 2025/2026"
$ws.Range("G52").Value = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
$ws.Range("G53").Value = "2025/2026"
$ws.Range("G54").Value = "2025/2026"
$ws.Range("G55").Value = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
$ws.Range("G56").Value = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
$ws.Range("G57").Value = "2025/2026"

# ---------------------------------------------------------------------
# 5) Class Statistics block (K2:L10) - recomputed totals now that two
#    more sessions have been recorded.
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 56      # Recorded Sessions
$ws.Range("L7").Value = 0       # Missing Sessions

# Percentage cells are stored as plain text in this report (matching the
# existing "96.4%"-style cells), so force text entry (otherwise Excel's
# input parser would coerce "100.0%" into a numeric percentage and swap
# in a numeric format) and then restore the original General-format look.
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "100.0%"   # Coverage %
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "43.5%"   # Average Attendance %

# ---------------------------------------------------------------------
# 6) Group Statistics block (K14:S16) - Year 3 / C1 and Year 3 / C2 rows.
# ---------------------------------------------------------------------
$ws.Range("O15").Value = 28
$ws.Range("P15").Value = 0
$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "100.0%"
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "48.5%"

$ws.Range("O16").Value = 28
$ws.Range("P16").Value = 0
$ws.Range("R16").NumberFormat = "@"
$ws.Range("R16").Value = "100.0%"
$ws.Range("S16").NumberFormat = "@"
$ws.Range("S16").Value = "38.5%"

# Restore the original (General-format) look on every percentage cell we
# just force-typed as text, by copying formatting back from an untouched
# neighbor that still carries the report's normal data-cell style.
$ws.Range("K4").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("R15").PasteSpecial(-4122)
$ws.Range("S15").PasteSpecial(-4122)
$ws.Range("R16").PasteSpecial(-4122)
$ws.Range("S16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
